# Apply the "units out of non-micro lab test" change:
# - Remove the unit portion "(D)" from the H-column label-building formulas
# - Add a new J column ("unit") that pulls the raw unit value from column D

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J
$ws.Cells.Item(1, 10).Value = "unit"
$ws.Cells.Item(1, 10).Style = $ws.Cells.Item(1, 9).Style

$firstRow = 3
$lastRow = 19

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Updated H formula - drop the " (" & D & ")" unit segment
    $ws.Cells.Item($r, 8).Formula = '=TRIM(UPPER(LEFT(C' + $r + ',1))&RIGHT(C' + $r + ',LEN(C' + $r + ')-1)&IF(E' + $r + '="","",", "&E' + $r + ')&IF(F' + $r + '="","",", at "&F' + $r + '))'

    # New J formula - surface the unit on its own
    $ws.Cells.Item($r, 10).Formula = '=IF(D' + $r + '="","",D' + $r + ')'
    $ws.Cells.Item($r, 10).Style = $ws.Cells.Item($r, 4).Style
}

# Selection / pane cosmetics to mirror the saved state
$ws.Range("B3").Select()
